$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 56 values
$ws.Range("B56").Value = 6.6
$ws.Range("D56").Value = 6.7

# Add new row 57
$ws.Range("A57").NumberFormat = "@"
$ws.Range("A57").Value = "01-08-2021"
$ws.Range("A57").Style = "Normal"
$ws.Range("B57").Value = 6.4
$ws.Range("C57").Value = 1.6
$ws.Range("D57").Value = 6.7
